# Adding the changes we made on may 9th
# Append 10 new data rows (22-31) of x/y/z gyroscope readings below the
# existing table in Sheet1 (which currently holds a header row + rows 2-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @(-0.0157297793775796, -0.0062613687478005, 0.1360702365636825),
    @(-0.0059559359215199, -0.0456621758639812, -0.009010262787342),
    @(0.0004581489483825,  -0.0007635815418325, -0.0609338097274303),
    @(-0.0006108652451075, 0.0114537235349416,  0.0355829000473022),
    @(-0.0042760567739605, -0.0050396383740007, -0.0058032199740409),
    @(0.005192354787141,   -0.0007635815418325, -0.0145080499351024),
    @(-0.0010690141934901, 0.008399397134780801, 0.0061086523346602),
    @(0.0004581489483825,  0.0038179077673703,  -0.0076358155347406),
    @(0.0007635815418325,  0.0027488935738801,  -0.0044287731871008),
    @(-0.00167987938039,   -0.0027488935738801, 0.0114537235349416)
)

$startRow = 22
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
